$wb = $excel.ActiveWorkbook

$updates = @(
    @{Sheet="ALC"; Cell="H2"; Value=109.416664},
    @{Sheet="ALC"; Cell="I2"; Value=113.3},
    @{Sheet="ALC"; Cell="K2"; Value=113.3},
    @{Sheet="ALC"; Cell="M2"; Value=-0.2999999999999972},
    @{Sheet="ALC"; Cell="H9"; Value=120.22222},
    @{Sheet="ALC"; Cell="I9"; Value=25.75},
    @{Sheet="ALC"; Cell="J9"; Value=195.8},
    @{Sheet="ALC"; Cell="K9"; Value=25.75},
    @{Sheet="ALC"; Cell="L9"; Value=195.8},
    @{Sheet="ALC"; Cell="M9"; Value=143.25},
    @{Sheet="ALC"; Cell="N9"; Value=-533.8},
    @{Sheet="ALC"; Cell="H17"; Value=528042.3},
    @{Sheet="ALC"; Cell="J17"; Value=557318.9399999999},
    @{Sheet="ALC"; Cell="L17"; Value=1671956.82},
    @{Sheet="ALC"; Cell="N17"; Value=-1672292.82},
    @{Sheet="ALC"; Cell="H55"; Value=798.2},
    @{Sheet="ALC"; Cell="I55"; Value=297.5},
    @{Sheet="ALC"; Cell="K55"; Value=297.5},
    @{Sheet="ALC"; Cell="M55"; Value=-83.5},
    @{Sheet="ALC"; Cell="H62"; Value=0},
    @{Sheet="ALC"; Cell="I62"; Value=0},
    @{Sheet="ALC"; Cell="J62"; Value=0},
    @{Sheet="ALC"; Cell="K62"; Value=0},
    @{Sheet="ALC"; Cell="L62"; Value=0},
    @{Sheet="ALC"; Cell="M62"; Value=$null},
    @{Sheet="ALC"; Cell="N62"; Value=$null},
    @{Sheet="ALC"; Cell="H65"; Value=0},
    @{Sheet="ALC"; Cell="I65"; Value=0},
    @{Sheet="ALC"; Cell="J65"; Value=0},
    @{Sheet="ALC"; Cell="K65"; Value=0},
    @{Sheet="ALC"; Cell="L65"; Value=0},
    @{Sheet="ALC"; Cell="M65"; Value=$null},
    @{Sheet="ALC"; Cell="N65"; Value=$null},
    @{Sheet="ALC"; Cell="H74"; Value=3217.9092},
    @{Sheet="ALC"; Cell="I74"; Value=2253},
    @{Sheet="ALC"; Cell="K74"; Value=2253},
    @{Sheet="ALC"; Cell="M74"; Value=-1317},
    @{Sheet="ALC"; Cell="H77"; Value=3217.9092},
    @{Sheet="ALC"; Cell="I77"; Value=2253},
    @{Sheet="ALC"; Cell="K77"; Value=11265},
    @{Sheet="ALC"; Cell="M77"; Value=-6585},
    @{Sheet="ALC"; Cell="H98"; Value=1998.875},
    @{Sheet="ALC"; Cell="I98"; Value=893.4737},
    @{Sheet="ALC"; Cell="K98"; Value=893.4737},
    @{Sheet="ALC"; Cell="M98"; Value=604.5263},
    @{Sheet="ALC"; Cell="H122"; Value=1998.875},
    @{Sheet="ALC"; Cell="I122"; Value=893.4737},
    @{Sheet="ALC"; Cell="K122"; Value=2680.4211},
    @{Sheet="ALC"; Cell="M122"; Value=-230.4211},
    @{Sheet="ARM"; Cell="H32"; Value=3029.8408},
    @{Sheet="ARM"; Cell="I32"; Value=2460.4},
    @{Sheet="ARM"; Cell="K32"; Value=2460.4},
    @{Sheet="ARM"; Cell="M32"; Value=-2173.4},
    @{Sheet="ARM"; Cell="H61"; Value=1731},
    @{Sheet="ARM"; Cell="I61"; Value=1418.2727},
    @{Sheet="ARM"; Cell="K61"; Value=1418.2727},
    @{Sheet="ARM"; Cell="M61"; Value=-1206.2727},
    @{Sheet="ARM"; Cell="H109"; Value=64000},
    @{Sheet="ARM"; Cell="J109"; Value=64000},
    @{Sheet="ARM"; Cell="L109"; Value=64000},
    @{Sheet="ARM"; Cell="N109"; Value=-66774},
    @{Sheet="ARM"; Cell="H110"; Value=1021.3684},
    @{Sheet="ARM"; Cell="I110"; Value=1083.5625},
    @{Sheet="ARM"; Cell="K110"; Value=1083.5625},
    @{Sheet="ARM"; Cell="M110"; Value=961.4375},
    @{Sheet="ARM"; Cell="H112"; Value=0},
    @{Sheet="ARM"; Cell="J112"; Value=0},
    @{Sheet="ARM"; Cell="L112"; Value=0},
    @{Sheet="ARM"; Cell="N112"; Value=$null},
    @{Sheet="ARM"; Cell="H118"; Value=64572.43},
    @{Sheet="ARM"; Cell="J118"; Value=64572.43},
    @{Sheet="ARM"; Cell="L118"; Value=64572.43},
    @{Sheet="ARM"; Cell="N118"; Value=-67886.42999999999},
    @{Sheet="ARM"; Cell="H132"; Value=1402.8478},
    @{Sheet="ARM"; Cell="I132"; Value=1251.8334},
    @{Sheet="ARM"; Cell="K132"; Value=3755.5002},
    @{Sheet="ARM"; Cell="M132"; Value=-1225.5002},
    @{Sheet="ARM"; Cell="H136"; Value=1731},
    @{Sheet="ARM"; Cell="I136"; Value=1418.2727},
    @{Sheet="ARM"; Cell="K136"; Value=4254.8181},
    @{Sheet="ARM"; Cell="M136"; Value=-1704.8181},
    @{Sheet="BSM"; Cell="H94"; Value=1167.1052},
    @{Sheet="BSM"; Cell="I94"; Value=1246.1765},
    @{Sheet="BSM"; Cell="K94"; Value=1246.1765},
    @{Sheet="BSM"; Cell="M94"; Value=-795.1765},
    @{Sheet="BSM"; Cell="H134"; Value=1410.2759},
    @{Sheet="BSM"; Cell="I134"; Value=1259.4},
    @{Sheet="BSM"; Cell="J134"; Value=2353.25},
    @{Sheet="BSM"; Cell="K134"; Value=3778.2},
    @{Sheet="BSM"; Cell="L134"; Value=7059.75},
    @{Sheet="BSM"; Cell="M134"; Value=-1243.2},
    @{Sheet="BSM"; Cell="N134"; Value=-12129.75},
    @{Sheet="CRP"; Cell="H58"; Value=1424.75},
    @{Sheet="CRP"; Cell="I58"; Value=1199.6666},
    @{Sheet="CRP"; Cell="J58"; Value=2100},
    @{Sheet="CRP"; Cell="K58"; Value=1199.6666},
    @{Sheet="CRP"; Cell="L58"; Value=2100},
    @{Sheet="CRP"; Cell="M58"; Value=-996.6666},
    @{Sheet="CRP"; Cell="N58"; Value=-2506},
    @{Sheet="CRP"; Cell="H94"; Value=2461.25},
    @{Sheet="CRP"; Cell="I94"; Value=2395.8},
    @{Sheet="CRP"; Cell="K94"; Value=2395.8},
    @{Sheet="CRP"; Cell="M94"; Value=-1944.8},
    @{Sheet="CRP"; Cell="H99"; Value=67499.75},
    @{Sheet="CRP"; Cell="I99"; Value=67499.75},
    @{Sheet="CRP"; Cell="K99"; Value=67499.75},
    @{Sheet="CRP"; Cell="M99"; Value=-66001.75},
    @{Sheet="CRP"; Cell="H122"; Value=116470},
    @{Sheet="CRP"; Cell="I122"; Value=200802},
    @{Sheet="CRP"; Cell="J122"; Value=11055},
    @{Sheet="CRP"; Cell="K122"; Value=602406},
    @{Sheet="CRP"; Cell="L122"; Value=33165},
    @{Sheet="CRP"; Cell="M122"; Value=-599956},
    @{Sheet="CRP"; Cell="N122"; Value=-38065},
    @{Sheet="CRP"; Cell="H126"; Value=67499.75},
    @{Sheet="CRP"; Cell="I126"; Value=67499.75},
    @{Sheet="CRP"; Cell="K126"; Value=202499.25},
    @{Sheet="CRP"; Cell="M126"; Value=-200029.25},
    @{Sheet="CRP"; Cell="H134"; Value=3956},
    @{Sheet="CRP"; Cell="I134"; Value=2305.8},
    @{Sheet="CRP"; Cell="J134"; Value=6313.4287},
    @{Sheet="CRP"; Cell="K134"; Value=6917.400000000001},
    @{Sheet="CRP"; Cell="L134"; Value=18940.2861},
    @{Sheet="CRP"; Cell="M134"; Value=-4382.400000000001},
    @{Sheet="CRP"; Cell="N134"; Value=-24010.2861},
    @{Sheet="CRP"; Cell="H136"; Value=1424.75},
    @{Sheet="CRP"; Cell="I136"; Value=1199.6666},
    @{Sheet="CRP"; Cell="J136"; Value=2100},
    @{Sheet="CRP"; Cell="K136"; Value=3598.9998},
    @{Sheet="CRP"; Cell="L136"; Value=6300},
    @{Sheet="CRP"; Cell="M136"; Value=-1048.9998},
    @{Sheet="CRP"; Cell="N136"; Value=-11400},
    @{Sheet="CUL"; Cell="H9"; Value=14094},
    @{Sheet="CUL"; Cell="I9"; Value=150},
    @{Sheet="CUL"; Cell="J9"; Value=35010},
    @{Sheet="CUL"; Cell="K9"; Value=450},
    @{Sheet="CUL"; Cell="L9"; Value=105030},
    @{Sheet="CUL"; Cell="M9"; Value=-226},
    @{Sheet="CUL"; Cell="N9"; Value=-105478},
    @{Sheet="CUL"; Cell="H10"; Value=169.4},
    @{Sheet="CUL"; Cell="I10"; Value=149},
    @{Sheet="CUL"; Cell="J10"; Value=251},
    @{Sheet="CUL"; Cell="K10"; Value=447},
    @{Sheet="CUL"; Cell="L10"; Value=753},
    @{Sheet="CUL"; Cell="M10"; Value=-308},
    @{Sheet="CUL"; Cell="N10"; Value=-1031},
    @{Sheet="CUL"; Cell="H13"; Value=78.75},
    @{Sheet="CUL"; Cell="I13"; Value=85.2},
    @{Sheet="CUL"; Cell="J13"; Value=68},
    @{Sheet="CUL"; Cell="K13"; Value=255.6},
    @{Sheet="CUL"; Cell="L13"; Value=204},
    @{Sheet="CUL"; Cell="M13"; Value=-87.60000000000002},
    @{Sheet="CUL"; Cell="N13"; Value=-540},
    @{Sheet="CUL"; Cell="H17"; Value=1000},
    @{Sheet="CUL"; Cell="I17"; Value=0},
    @{Sheet="CUL"; Cell="J17"; Value=1000},
    @{Sheet="CUL"; Cell="K17"; Value=0},
    @{Sheet="CUL"; Cell="L17"; Value=3000},
    @{Sheet="CUL"; Cell="M17"; Value=$null},
    @{Sheet="CUL"; Cell="N17"; Value=-3338},
    @{Sheet="GSM"; Cell="H80"; Value=2789},
    @{Sheet="GSM"; Cell="I80"; Value=2083.3333},
    @{Sheet="GSM"; Cell="J80"; Value=3847.5},
    @{Sheet="GSM"; Cell="K80"; Value=2083.3333},
    @{Sheet="GSM"; Cell="L80"; Value=3847.5},
    @{Sheet="GSM"; Cell="M80"; Value=-1085.3333},
    @{Sheet="GSM"; Cell="N80"; Value=-5843.5},
    @{Sheet="GSM"; Cell="H83"; Value=2789},
    @{Sheet="GSM"; Cell="I83"; Value=2083.3333},
    @{Sheet="GSM"; Cell="J83"; Value=3847.5},
    @{Sheet="GSM"; Cell="K83"; Value=10416.6665},
    @{Sheet="GSM"; Cell="L83"; Value=19237.5},
    @{Sheet="GSM"; Cell="M83"; Value=-5424.666499999999},
    @{Sheet="GSM"; Cell="N83"; Value=-29221.5},
    @{Sheet="GSM"; Cell="H122"; Value=4196.75},
    @{Sheet="GSM"; Cell="J122"; Value=5000},
    @{Sheet="GSM"; Cell="L122"; Value=15000},
    @{Sheet="GSM"; Cell="N122"; Value=-19900},
    @{Sheet="LTW"; Cell="H110"; Value=0},
    @{Sheet="LTW"; Cell="J110"; Value=0},
    @{Sheet="LTW"; Cell="L110"; Value=0},
    @{Sheet="LTW"; Cell="N110"; Value=$null},
    @{Sheet="LTW"; Cell="H118"; Value=36000},
    @{Sheet="LTW"; Cell="J118"; Value=36000},
    @{Sheet="LTW"; Cell="L118"; Value=36000},
    @{Sheet="LTW"; Cell="N118"; Value=-39314},
    @{Sheet="WVR"; Cell="H42"; Value=50000},
    @{Sheet="WVR"; Cell="I42"; Value=50000},
    @{Sheet="WVR"; Cell="K42"; Value=50000},
    @{Sheet="WVR"; Cell="M42"; Value=-49622},
    @{Sheet="WVR"; Cell="H43"; Value=29027},
    @{Sheet="WVR"; Cell="I43"; Value=29027},
    @{Sheet="WVR"; Cell="K43"; Value=29027},
    @{Sheet="WVR"; Cell="M43"; Value=-28878},
    @{Sheet="WVR"; Cell="H107"; Value=18520160},
    @{Sheet="WVR"; Cell="I107"; Value=970.4375},
    @{Sheet="WVR"; Cell="K107"; Value=2911.3125},
    @{Sheet="WVR"; Cell="M107"; Value=-991.3125},
    @{Sheet="WVR"; Cell="H126"; Value=1873.875},
    @{Sheet="WVR"; Cell="I126"; Value=1761},
    @{Sheet="WVR"; Cell="K126"; Value=5283},
    @{Sheet="WVR"; Cell="M126"; Value=-2813}
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    if ($u.Value -eq $null) {
        $ws.Range($u.Cell).ClearContents()
    } else {
        $ws.Range($u.Cell).Value = $u.Value
    }
}

Write-Output "Applied $($updates.Count) cell updates."